# Edit script: insert 3 new data rows (new weekly price records, date 44769)
# into the "Hortaliza, Femacal de La Calera - Tomate" sheet, right before the
# existing row that used to be row 1145 (date 44596 / Primera). All rows from
# the old row 1145 through 1238 shift down by 3 rows (to 1148..1241) and keep
# their original values unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 1145:1147, pushing the old 1145..1238 block down to 1148..1241.
$ws.Rows("1145:1147").Insert()

# Populate the 3 newly inserted rows with the new weekly records.

# Row 1145
$ws.Range("A1145").Value = 3
$ws.Range("B1145").Value = "Femacal de La Calera"
$ws.Range("C1145").Value = "Coquimbo"
$ws.Range("D1145").Value = 44769
$ws.Range("D1145").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E1145").Value = 5
$ws.Range("F1145").Value = 100112020
$ws.Range("G1145").Value = "Tomate"
$ws.Range("H1145").Value = "Larga vida"
$ws.Range("I1145").Value = "Primera"
$ws.Range("J1145").Value = 510
$ws.Range("K1145").Value = 9000
$ws.Range("L1145").Value = 9500
$ws.Range("M1145").Value = 9255
$ws.Range("N1145").Value = "`$/bandeja 18 kilos"
$ws.Range("O1145").Value = "Región de Arica y Parinacota"
$ws.Range("P1145").Value = 514
$ws.Range("Q1145").Value = 18
$ws.Range("R1145").Value = "Hortaliza"

# Row 1146
$ws.Range("A1146").Value = 3
$ws.Range("B1146").Value = "Femacal de La Calera"
$ws.Range("C1146").Value = "Coquimbo"
$ws.Range("D1146").Value = 44769
$ws.Range("D1146").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E1146").Value = 5
$ws.Range("F1146").Value = 100112020
$ws.Range("G1146").Value = "Tomate"
$ws.Range("H1146").Value = "Larga vida"
$ws.Range("I1146").Value = "Primera"
$ws.Range("J1146").Value = 570
$ws.Range("K1146").Value = 3800
$ws.Range("L1146").Value = 4000
$ws.Range("M1146").Value = 3923
$ws.Range("N1146").Value = "`$/caja 12 kilos"
$ws.Range("O1146").Value = "Región de Arica y Parinacota"
$ws.Range("P1146").Value = 327
$ws.Range("Q1146").Value = 12
$ws.Range("R1146").Value = "Hortaliza"

# Row 1147
$ws.Range("A1147").Value = 3
$ws.Range("B1147").Value = "Femacal de La Calera"
$ws.Range("C1147").Value = "Coquimbo"
$ws.Range("D1147").Value = 44769
$ws.Range("D1147").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E1147").Value = 5
$ws.Range("F1147").Value = 100112020
$ws.Range("G1147").Value = "Tomate"
$ws.Range("H1147").Value = "Larga vida"
$ws.Range("I1147").Value = "Segunda"
$ws.Range("J1147").Value = 280
$ws.Range("K1147").Value = 7000
$ws.Range("L1147").Value = 7000
$ws.Range("M1147").Value = 7000
$ws.Range("N1147").Value = "`$/bandeja 18 kilos"
$ws.Range("O1147").Value = "Región de Arica y Parinacota"
$ws.Range("P1147").Value = 389
$ws.Range("Q1147").Value = 18
$ws.Range("R1147").Value = "Hortaliza"
